# Generated BOM and CPL files
# Row 6 (R3,R2,R8,R9,R5,R4 @ qty 6) loses designator "R2": qty drops to 5
# and the designator list becomes "R3,R8,R9,R5,R4".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "R3,R8,R9,R5,R4"

# Match the author's final cursor position recorded in the sheet view.
[void]$ws.Range("B7").Select()
